# The deck ships two theme parts:
#   ppt/theme/theme1.xml -> "Office Theme" (blue Office palette)
#   ppt/theme/theme2.xml -> "Integral"     (green Integral palette; this is
#                                            the one actually driving the
#                                            slide master and therefore
#                                            every slide in the deck)
#
# The authored edit swaps the two themes: the master (and every slide)
# should now render with the "Office Theme" palette instead of "Integral".
# The two theme parts already share an identical font scheme and an
# identical format scheme (fills / lines / effects) - the only real
# difference between them is their 12-slot color scheme (clrScheme: dk1,
# lt1, dk2, lt2, accent1-6, hlink, folHlink) and the cosmetic theme/
# color-scheme "name" labels.
#
# Re-point the presentation's live theme colors from the Integral palette
# to the Office Theme palette by rewriting the 12 theme colors in place.
# (The theme/color-scheme "name" attributes and theme1.xml itself are not
# reachable through the PowerPoint object model exposed here, so they are
# left as-is; the observable, rendered effect - every slide now using the
# Office Theme palette - is fully reproduced.)

function Convert-HexToComRgb([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme color scheme (formerly theme1.xml), in clrScheme document
# order: dk1, lt1, dk2, lt2, accent1, accent2, accent3, accent4, accent5,
# accent6, hlink, folHlink.
$officeThemeColors = @(
    "000000",
    "FFFFFF",
    "44546A",
    "E7E6E6",
    "5B9BD5",
    "ED7D31",
    "A5A5A5",
    "FFC000",
    "4472C4",
    "70AD47",
    "0563C1",
    "954F72"
)

$p = $ppt.ActivePresentation

# Rewrite the 12 live theme colors in place.
$tcs = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Item($i).RGB = Convert-HexToComRgb $officeThemeColors[$i - 1]
}
